$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Four cluster names used for the Sending cluster (col A) / Target cluster (col D) values.
$clusters = @("ECs", "FAPs", "M2", "sCs")

# Numeric payload (columns E..T) for worksheet rows 2..17, in row order.
# Row order follows the nested loop: outer = Sending cluster, inner = Target cluster,
# i.e. row 2..5 = ECs x {ECs,FAPs,M2,sCs}, row 6..9 = FAPs x {...}, etc.

$data = @(
    @(3, 1, 16.506869, 49.520607, 0.2165594803671733, 0.2165594803671733, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 285.9802398983767, 2573.82215908539, 0.03922692237743623, 0.03922692237743623),
    @(3, 1, 16.506869, 49.520607, 0.2165594803671733, 0.2165594803671733, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 108.7881845048303, 979.093660543473, 0.01492209976000308, 0.01492209976000308),
    @(3, 1, 16.506869, 49.520607, 0.2165594803671733, 0.2165594803671733, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 1.835481298455, 16.519331686095, 0.000251766634104914, 0.0002517666341049141),
    @(3, 1, 16.506869, 49.520607, 0.2165594803671733, 0.2165594803671733, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 1182.202903350884, 10639.82613015796, 0.1621586915956291, 0.1621586915956291),
    @(3, 1, 24.781512, 74.344536, 0.3251174623990092, 0.3251174623990092, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 429.3377954840801, 3864.04015935672, 0.05889078344412284, 0.05889078344412284),
    @(3, 1, 24.781512, 74.344536, 0.3251174623990092, 0.3251174623990092, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 163.322050945256, 1469.898458507304, 0.02240232198290988, 0.02240232198290988),
    @(3, 1, 24.781512, 74.344536, 0.3251174623990092, 0.3251174623990092, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 2.75558022684, 24.80022204156, 0.000377973428169239, 0.000377973428169239),
    @(3, 1, 24.781512, 74.344536, 0.3251174623990092, 0.3251174623990092, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 1774.823283314648, 15973.40954983183, 0.2434463835438073, 0.2434463835438072),
    @(3, 1, 24.67943933333333, 74.038318, 0.3237783348120013, 0.3237783348120013, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 427.5693943596511, 3848.12454923686, 0.05864821796594577, 0.05864821796594577),
    @(3, 1, 24.67943933333333, 74.038318, 0.3237783348120013, 0.3237783348120013, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 162.6493431110669, 1463.844087999602, 0.02231004897130667, 0.02231004897130667),
    @(3, 1, 24.67943933333333, 74.038318, 0.3237783348120013, 0.3237783348120013, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 2.74423025667, 24.69807231003, 0.0003764165919381658, 0.0003764165919381658),
    @(3, 1, 24.67943933333333, 74.038318, 0.3237783348120013, 0.3237783348120013, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 1767.512956753863, 15907.61661078477, 0.2424436512828107, 0.2424436512828107),
    @(3, 1, 10.255437, 30.766311, 0.1345447224218162, 0.1345447224218162, 3, 1, 17.32492333333333, 51.97477, 0.1811369435820938, 0.1811369435820938, 177.67465977483, 1599.07193797347, 0.024371019794589, 0.024371019794589),
    @(3, 1, 10.255437, 30.766311, 0.1345447224218162, 0.1345447224218162, 3, 1, 6.590479666666667, 19.771439, 0.06890531753540824, 0.06890531753540824, 67.588249021281, 608.2942411915291, 0.009270846821188606, 0.009270846821188608),
    @(3, 1, 10.255437, 30.766311, 0.1345447224218162, 0.1345447224218162, 2, 0.6666666666666666, 0.111195, 0.333585, 0.001162574982531578, 0.001162574982531578, 1.140353317215, 10.263179854935, 0.000156418328319259, 0.000156418328319259),
    @(3, 1, 10.255437, 30.766311, 0.1345447224218162, 0.1345447224218162, 3, 1, 71.61884566666667, 214.856537, 0.7487951638999665, 0.7487951638999664, 734.4825597472231, 6610.343037725008, 0.1007464374777194, 0.1007464374777194),
)

# Fill in rows 2..17: columns A (Sending cluster) and D (Target cluster) follow a
# 4x4 nested-loop pattern over $clusters; columns B/C are always "Psen1"/"Notch3"
# (left untouched, already correct); columns E..T come from $data.
$row = 2
for ($i = 0; $i -lt 4; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $ws.Cells.Item($row, 1).Value = $clusters[$i]
        $ws.Cells.Item($row, 2).Value = "Psen1"
        $ws.Cells.Item($row, 3).Value = "Notch3"
        $ws.Cells.Item($row, 4).Value = $clusters[$j]

        $values = $data[$row - 2]
        for ($c = 0; $c -lt $values.Length; $c++) {
            $ws.Cells.Item($row, 5 + $c).Value = $values[$c]
        }

        $row++
    }
}
